$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 175045
$ws.Range("C4").Value = 165032
$ws.Range("C7").Value = 5.72
$ws.Range("C8").Value = 64.56999999999999
